$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "ano" in A2
$ws.Range("A2").Value = "ano"

# Fill A3:A14 with the year value 2023
$ws.Range("A3:A14").Value = 2023

# Update the active selection to A2 (matches new sheetView selection)
$ws.Range("A2").Select()
